# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "27.185.95"
    "E2"  = "  +0.56%  "
    "D3"  = "1.685.92"
    "E3"  = "  +0.24%  "
    "D5"  = "215.96"
    "E6"  = "  +0.62%  "
    "E7"  = "  +0.12%  "
    "D8"  = "23.09"
    "E8"  = "  +7.94%  "
    "E9"  = "  +3.41%  "
    "E10" = "  +0.92%  "
    "D11" = "0.0890"
    "E11" = "  +0.48%  "
    "D12" = "1.923.09"
    "E12" = "  +0.20%  "
    "D13" = "1.690.79"
    "E13" = "  -0.17%  "
    "E14" = "  +2.26%  "
    "E15" = "  +4.01%  "
    "D16" = "66.90"
    "E16" = "  +1.21%  "
    "D17" = "27.184.07"
    "E17" = "  +0.45%  "
    "D18" = "236.13"
    "E18" = "  -0.10%  "
    "D19" = "8.02"
    "E19" = "  -2.14%  "
    "D20" = "0.0₃0744"
    "E20" = "  +1.29%  "
    "E21" = "  +0.16%  "
    "E22" = "  +2.25%  "
    "D23" = "9.60"
    "E23" = "  +4.00%  "
    "E24" = "  -2.89%  "
    "D25" = "147.26"
    "E25" = "  +0.28%  "
    "D26" = "7.33"
    "E26" = "  +1.37%  "
    "D27" = "16.44"
    "E27" = "  +2.41%  "
    "E28" = "  +0.58%  "
    "E29" = "  +0.19%  "
    "D30" = "0.0506"
    "E30" = "  +1.16%  "
    "E31" = "  +0.13%  "
    "E32" = "  +1.27%  "
    "D33" = "1.545.94"
    "E33" = "  +2.21%  "
    "E34" = "  +1.39%  "
    "D35" = "1.67"
    "E35" = "  -1.12%  "
    "D36" = "0.605"
    "E36" = "  +2.78%  "
    "E37" = "  +3.06%  "
    "E38" = "  -0.44%  "
    "E39" = "  -0.37%  "
    "E40" = "  +1.55%  "
    "E41" = "  +1.10%  "
    "E43" = "  +0.13%  "
    "E44" = "  -0.87%  "
    "D45" = "1.831.67"
    "E45" = "  +0.46%  "
    "E46" = "  +1.06%  "
    "D47" = "90.22"
    "E47" = "  +0.15%  "
    "E48" = "  +5.22%  "
    "D49" = "1.62"
    "E49" = "  +5.91%  "
    "D50" = "8.30"
    "E50" = "  +5.63%  "
    "E51" = "  -0.55%  "
}

# All of these cells hold plain text (numbers-as-strings like prices, and
# percentage strings padded with spaces). Force the cell to stay text
# before writing so Excel doesn't auto-convert numeric-looking values
# (e.g. "215.96") into real numbers, then clear the formatting that the
# text coercion leaves behind so the cell's style stays the same as before.
foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
